$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style) of the last existing row (238) down to the
# new rows 239:244 so the date column keeps its style (s="2").
$ws.Range("A238").Copy()
$ws.Range("A239:A244").PasteSpecial(-4122)

# New data rows (dates are Excel serial date numbers).
$dates = @(44313, 44314, 44315, 44316, 44317, 44318)
$bvals = @(2, 0, 0, 1, 2, 0)
$cvals = @(2, 2, 2, 3, 5, 5)
$dvals = @(96.15384615384616, 96.15384615384616, 96.15384615384616, 144.2307692307692, 240.3846153846154, 240.3846153846154)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = 239 + $i
    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 2).Value = $bvals[$i]
    $ws.Cells.Item($r, 3).Value = $cvals[$i]
    $ws.Cells.Item($r, 4).Value = $dvals[$i]
}
